# Update the "Taches" planning sheet: condense the per-class task rows
# (8-16) down to two summary rows, clearing the rest.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taches")

# Row 8: new summary task covering every class
$ws.Range("A8").Value = "Ensemble des classes"
$ws.Range("B8").Value = "implémentation"
$ws.Range("C8").Value = "création du fichier de chaque classe"
# D8, E8, F8, G8, H8 stay as-is (Classe / 1 / Dorian / 1 / Mardi 08)

# Row 9: irrigation task, everything past column B cleared
$ws.Range("A9").Value = "Irrigation"
$ws.Range("B9").Value = "arrosageParcelles()"
$ws.Range("C9:H9").ClearContents()

# Rows 10-16: fully cleared (keep the row/cell styling, drop the values)
$ws.Range("A10:H16").ClearContents()

# Move the active-cell selection to C9, matching the saved view state
$ws.Range("C9").Select()
